# Status Summary Source File - apply edit described by commit:
# "Remove original Status Summary section, update revenue chart to show
#  revenue+expenses with totals, fix margin chart hover to show percentages"
#
# The portion of that change that lives in this workbook's own sheet
# (Sheet1) is:
#   - add a "Measures" header label in A1 (styled like the other row labels)
#   - drop the stale active-cell selection left over in the sheet view
#   - row 3 goes back to the sheet's default (auto) row height
#   - the 2024 Total Expenses figure for Mar (D7) is corrected to 124364,
#     which ripples into the 2024 Margin (D5) formula result automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Measures" label to A1, formatted like the other row headers (A2:A7) ---
$ws.Range("A1").Value = "Measures"
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 3 reverts to the sheet's default auto row height (no explicit ht) ---
$ws.Rows.Item(3).AutoFit()

# --- Correct 2024 Total Expenses for Mar (D7); D5 (2024 Margin) recalculates ---
$ws.Range("D7").Value = 124364

# --- Clear the stray selection left at I10, resetting it to the top-left cell ---
$ws.Range("A1").Select()
